$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-23 09:18:06"
$wsZh.Range("H2").Value = "2016-03-23 09:18:37"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-23 09:18:10"
$wsDe.Range("H2").Value = "2016-03-23 09:18:43"
